$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 31, pushing existing rows 31-100 down to 33-102.
$ws.Rows.Item(31).Resize(2).Insert()

# New row 31: Early Treat / Primera
$ws.Cells.Item(31, 1).Value = 1
$ws.Cells.Item(31, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(31, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(31, 4).Value = 45272
$ws.Cells.Item(31, 5).Value = 15
$ws.Cells.Item(31, 6).Value = "Fruta"
$ws.Cells.Item(31, 7).Value = 100103
$ws.Cells.Item(31, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(31, 9).Value = 100103004
$ws.Cells.Item(31, 10).Value = "Durazno"
$ws.Cells.Item(31, 11).Value = "Early Treat"
$ws.Cells.Item(31, 12).Value = "Primera"
$ws.Cells.Item(31, 13).Value = 300
$ws.Cells.Item(31, 14).Value = 24000
$ws.Cells.Item(31, 15).Value = 25000
$ws.Cells.Item(31, 16).Value = 24500
$ws.Cells.Item(31, 17).Value = "$/bandeja 18 kilos granel"
$ws.Cells.Item(31, 18).Value = "Provincia de San Felipe de Aconcagua"
$ws.Cells.Item(31, 19).Value = 1361
$ws.Cells.Item(31, 20).Value = 18

# New row 32: Florida King / Tercera
$ws.Cells.Item(32, 1).Value = 1
$ws.Cells.Item(32, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(32, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(32, 4).Value = 45272
$ws.Cells.Item(32, 5).Value = 15
$ws.Cells.Item(32, 6).Value = "Fruta"
$ws.Cells.Item(32, 7).Value = 100103
$ws.Cells.Item(32, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(32, 9).Value = 100103004
$ws.Cells.Item(32, 10).Value = "Durazno"
$ws.Cells.Item(32, 11).Value = "Florida King"
$ws.Cells.Item(32, 12).Value = "Tercera"
$ws.Cells.Item(32, 13).Value = 300
$ws.Cells.Item(32, 14).Value = 14000
$ws.Cells.Item(32, 15).Value = 15000
$ws.Cells.Item(32, 16).Value = 14500
$ws.Cells.Item(32, 17).Value = "$/bandeja 18 kilos granel"
$ws.Cells.Item(32, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(32, 19).Value = 806
$ws.Cells.Item(32, 20).Value = 18
